# Update automàtic: dades i banners [2026-02-07 03:49]
# Applies the meteocat daily-summary refresh: new DATA_EXTRACCIO timestamps
# and refreshed observation values (humidity, pressure, wind, temperature).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("E2").Value = '2026-02-07 03:47:40'
$ws.Range("H2").Value = "'94%"
$ws.Range("N2").Value = '-1.5 °C 3:29 TU'

# Row 3
$ws.Range("E3").Value = '2026-02-07 03:47:42'
$ws.Range("N3").Value = '-7.0 °C 3:01 TU'
$ws.Range("O3").Value = '-5.4 °C'

# Row 4
$ws.Range("E4").Value = '2026-02-07 03:47:45'
$ws.Range("H4").Value = "'56%"
$ws.Range("J4").Value = '1000.8 hPa'
$ws.Range("O4").Value = '11.8 °C'

# Row 5
$ws.Range("E5").Value = '2026-02-07 03:47:47'
$ws.Range("H5").Value = "'72%"
$ws.Range("J5").Value = '1000.8 hPa'
$ws.Range("O5").Value = '9.0 °C'

# Row 6
$ws.Range("E6").Value = '2026-02-07 03:47:49'
$ws.Range("J6").Value = '1002.4 hPa'
$ws.Range("N6").Value = '11.4 °C 3:16 TU'
$ws.Range("O6").Value = '12.0 °C'

# Row 7
$ws.Range("E7").Value = '2026-02-07 03:47:52'
$ws.Range("J7").Value = '1002.3 hPa'
$ws.Range("N7").Value = '7.1 °C 3:23 TU'
$ws.Range("O7").Value = '8.0 °C'

# Row 8
$ws.Range("E8").Value = '2026-02-07 03:47:54'
$ws.Range("N8").Value = '3.1 °C 3:25 TU'
$ws.Range("O8").Value = '4.7 °C'

# Row 9
$ws.Range("E9").Value = '2026-02-07 03:47:56'
$ws.Range("N9").Value = '1.3 °C 3:28 TU'
$ws.Range("O9").Value = '2.4 °C'

# Row 10
$ws.Range("E10").Value = '2026-02-07 03:47:59'
$ws.Range("N10").Value = '6.2 °C 3:18 TU'
$ws.Range("O10").Value = '7.1 °C'

# Row 11
$ws.Range("E11").Value = '2026-02-07 03:48:01'
$ws.Range("O11").Value = '1.2 °C'

# Row 12
$ws.Range("E12").Value = '2026-02-07 03:48:03'

# Row 13
$ws.Range("E13").Value = '2026-02-07 03:48:05'
$ws.Range("H13").Value = "'91%"
$ws.Range("O13").Value = '7.1 °C'

# Row 14
$ws.Range("E14").Value = '2026-02-07 03:48:08'
$ws.Range("H14").Value = "'81%"
$ws.Range("O14").Value = '-5.6 °C'

# Row 15
$ws.Range("E15").Value = '2026-02-07 03:48:10'
$ws.Range("H15").Value = "'79%"
$ws.Range("N15").Value = '4.9 °C 3:29 TU'
$ws.Range("O15").Value = '7.6 °C'

# Row 16
$ws.Range("E16").Value = '2026-02-07 03:48:12'
$ws.Range("H16").Value = "'88%"
$ws.Range("N16").Value = '2.1 °C 3:28 TU'
$ws.Range("O16").Value = '3.4 °C'

# Row 17
$ws.Range("E17").Value = '2026-02-07 03:48:15'
$ws.Range("J17").Value = '1004.4 hPa'
$ws.Range("N17").Value = '2.9 °C 3:00 TU'
$ws.Range("O17").Value = '3.5 °C'

# Row 18
$ws.Range("E18").Value = '2026-02-07 03:48:17'
$ws.Range("N18").Value = '-8.0 °C 3:29 TU'
$ws.Range("O18").Value = '-6.7 °C'

# Row 19
$ws.Range("E19").Value = '2026-02-07 03:48:20'
$ws.Range("J19").Value = '1005.5 hPa'
$ws.Range("N19").Value = '4.1 °C 3:17 TU'
$ws.Range("O19").Value = '4.9 °C'

# Row 20
$ws.Range("E20").Value = '2026-02-07 03:48:22'
$ws.Range("H20").Value = "'88%"

# Row 21
$ws.Range("E21").Value = '2026-02-07 03:48:25'
$ws.Range("H21").Value = "'68%"
$ws.Range("J21").Value = '1001.0 hPa'
$ws.Range("N21").Value = '4.7 °C 3:27 TU'
$ws.Range("O21").Value = '8.3 °C'

# Row 22
$ws.Range("E22").Value = '2026-02-07 03:48:27'
$ws.Range("N22").Value = '4.3 °C 3:27 TU'
$ws.Range("O22").Value = '5.7 °C'

# Row 23
$ws.Range("E23").Value = '2026-02-07 03:48:30'
$ws.Range("J23").Value = '1001.0 hPa'
$ws.Range("N23").Value = '7.1 °C 3:11 TU'
$ws.Range("O23").Value = '7.7 °C'

# Row 24
$ws.Range("E24").Value = '2026-02-07 03:48:32'
$ws.Range("H24").Value = "'80%"
$ws.Range("J24").Value = '1000.4 hPa'
$ws.Range("L24").Value = '27.0 km/h - 357º 3:21 TU'
$ws.Range("N24").Value = '9.9 °C 3:14 TU'

# Row 25
$ws.Range("E25").Value = '2026-02-07 03:48:35'
$ws.Range("O25").Value = '0.7 °C'

# Row 26
$ws.Range("E26").Value = '2026-02-07 03:48:37'
$ws.Range("H26").Value = "'75%"
$ws.Range("L26").Value = '35.3 km/h - 44º 3:09 TU'

# Row 27
$ws.Range("E27").Value = '2026-02-07 03:48:39'
$ws.Range("H27").Value = "'96%"
$ws.Range("N27").Value = '7.6 °C 3:16 TU'
$ws.Range("O27").Value = '8.4 °C'

# Row 28
$ws.Range("E28").Value = '2026-02-07 03:48:42'
$ws.Range("J28").Value = '1003.2 hPa'
$ws.Range("N28").Value = '2.5 °C 3:19 TU'
$ws.Range("O28").Value = '3.8 °C'

# Row 29
$ws.Range("E29").Value = '2026-02-07 03:48:44'
$ws.Range("N29").Value = '10.4 °C 3:25 TU'
$ws.Range("O29").Value = '11.6 °C'

# Row 30
$ws.Range("E30").Value = '2026-02-07 03:48:46'
$ws.Range("H30").Value = "'84%"
$ws.Range("L30").Value = '34.2 km/h - 327º 3:18 TU'
$ws.Range("O30").Value = '-4.7 °C'

# Row 31
$ws.Range("E31").Value = '2026-02-07 03:48:49'

# Row 32
$ws.Range("E32").Value = '2026-02-07 03:48:51'
$ws.Range("H32").Value = "'62%"
$ws.Range("J32").Value = '1003.8 hPa'
$ws.Range("O32").Value = '11.4 °C'

# Row 33
$ws.Range("E33").Value = '2026-02-07 03:48:53'
$ws.Range("N33").Value = '6.5 °C 3:27 TU'
$ws.Range("O33").Value = '7.6 °C'

# Row 34
$ws.Range("E34").Value = '2026-02-07 03:48:56'
$ws.Range("N34").Value = '5.4 °C 3:08 TU'
$ws.Range("O34").Value = '6.8 °C'

# Row 35
$ws.Range("E35").Value = '2026-02-07 03:48:58'
$ws.Range("N35").Value = '-6.3 °C 3:29 TU'
$ws.Range("O35").Value = '-4.4 °C'

# Row 36
$ws.Range("E36").Value = '2026-02-07 03:49:00'
$ws.Range("J36").Value = '1006.0 hPa'
$ws.Range("N36").Value = '4.1 °C 3:18 TU'
